$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 909.2353000000001
$ws.Range("J19").Value = 920.1111
$ws.Range("L19").Value = 920.1111
$ws.Range("N19").Value = -1270.1111
$ws.Range("H98").Value = 1236.6
$ws.Range("I98").Value = 858.25
$ws.Range("K98").Value = 858.25
$ws.Range("M98").Value = 639.75
$ws.Range("H106").Value = 372932.06
$ws.Range("I106").Value = 464448.5
$ws.Range("K106").Value = 464448.5
$ws.Range("M106").Value = -463817.5
$ws.Range("H109").Value = 88975
$ws.Range("J109").Value = 88975
$ws.Range("L109").Value = 88975
$ws.Range("N109").Value = -91749
$ws.Range("H122").Value = 1236.6
$ws.Range("I122").Value = 858.25
$ws.Range("K122").Value = 2574.75
$ws.Range("M122").Value = -124.75
$ws.Range("H132").Value = 1744.7567
$ws.Range("I132").Value = 1274.6774
$ws.Range("K132").Value = 3824.0322
$ws.Range("M132").Value = -1294.0322
$ws.Range("H137").Value = 660873.5600000001
$ws.Range("I137").Value = 1511.5385
$ws.Range("K137").Value = 4534.6155
$ws.Range("M137").Value = -1984.6155
$ws.Range("H138").Value = 2203.2126
$ws.Range("I138").Value = 863.55
$ws.Range("K138").Value = 2590.65
$ws.Range("M138").Value = 2549.35

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 592.8158
$ws.Range("I2").Value = 457.48486
$ws.Range("J2").Value = 1486
$ws.Range("K2").Value = 457.48486
$ws.Range("L2").Value = 1486
$ws.Range("M2").Value = -344.48486
$ws.Range("N2").Value = -1712
$ws.Range("H4").Value = 451.55554
$ws.Range("J4").Value = 417.5
$ws.Range("L4").Value = 417.5
$ws.Range("N4").Value = -649.5
$ws.Range("H41").Value = 3226.6667
$ws.Range("I41").Value = 3226.6667
$ws.Range("K41").Value = 3226.6667
$ws.Range("M41").Value = -2812.6667
$ws.Range("H74").Value = 2429.1924
$ws.Range("I74").Value = 1669.3125
$ws.Range("K74").Value = 1669.3125
$ws.Range("M74").Value = -795.3125
$ws.Range("H77").Value = 2429.1924
$ws.Range("I77").Value = 1669.3125
$ws.Range("K77").Value = 8346.5625
$ws.Range("M77").Value = -3978.5625
$ws.Range("H102").Value = 45757.92
$ws.Range("I102").Value = 47052.047
$ws.Range("K102").Value = 47052.047
$ws.Range("M102").Value = -45430.047
$ws.Range("H116").Value = 592.8158
$ws.Range("I116").Value = 457.48486
$ws.Range("J116").Value = 1486
$ws.Range("K116").Value = 457.48486
$ws.Range("L116").Value = 1486
$ws.Range("M116").Value = 1836.51514
$ws.Range("N116").Value = -6074
$ws.Range("H125").Value = 54775.668
$ws.Range("J125").Value = 54775.668
$ws.Range("L125").Value = 54775.668
$ws.Range("N125").Value = -64615.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 592.8158
$ws.Range("I3").Value = 457.48486
$ws.Range("J3").Value = 1486
$ws.Range("K3").Value = 457.48486
$ws.Range("L3").Value = 1486
$ws.Range("M3").Value = -343.48486
$ws.Range("N3").Value = -1714
$ws.Range("H86").Value = 3708.3
$ws.Range("I86").Value = 2678.3333
$ws.Range("J86").Value = 5253.25
$ws.Range("K86").Value = 2678.3333
$ws.Range("L86").Value = 5253.25
$ws.Range("M86").Value = -1555.3333
$ws.Range("N86").Value = -7499.25
$ws.Range("H89").Value = 3708.3
$ws.Range("I89").Value = 2678.3333
$ws.Range("J89").Value = 5253.25
$ws.Range("K89").Value = 13391.6665
$ws.Range("L89").Value = 26266.25
$ws.Range("M89").Value = -7775.666499999999
$ws.Range("N89").Value = -37498.25
$ws.Range("H94").Value = 640.5625
$ws.Range("I94").Value = 599.2727
$ws.Range("J94").Value = 731.4
$ws.Range("K94").Value = 599.2727
$ws.Range("L94").Value = 731.4
$ws.Range("M94").Value = -148.2727
$ws.Range("N94").Value = -1633.4
$ws.Range("H112").Value = 99988.336
$ws.Range("J112").Value = 99988.336
$ws.Range("L112").Value = 99988.336
$ws.Range("N112").Value = -102942.336
$ws.Range("H134").Value = 1775
$ws.Range("I134").Value = 1354.6666
$ws.Range("K134").Value = 4063.9998
$ws.Range("M134").Value = -1528.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2020.2142
$ws.Range("I58").Value = 1836.0834
$ws.Range("J58").Value = 3125
$ws.Range("K58").Value = 1836.0834
$ws.Range("L58").Value = 3125
$ws.Range("M58").Value = -1633.0834
$ws.Range("N58").Value = -3531
$ws.Range("H105").Value = 64438.777
$ws.Range("J105").Value = 4181.6
$ws.Range("L105").Value = 4181.6
$ws.Range("N105").Value = -7675.6
$ws.Range("H136").Value = 2020.2142
$ws.Range("I136").Value = 1836.0834
$ws.Range("J136").Value = 3125
$ws.Range("K136").Value = 5508.2502
$ws.Range("L136").Value = 9375
$ws.Range("M136").Value = -2958.2502
$ws.Range("N136").Value = -14475
$ws.Range("H141").Value = 215606.19
$ws.Range("J141").Value = 215606.19
$ws.Range("L141").Value = 215606.19
$ws.Range("N141").Value = -225966.19

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 395
$ws.Range("I63").Value = 542.5
$ws.Range("K63").Value = 1627.5
$ws.Range("M63").Value = -878.5
$ws.Range("H66").Value = 395
$ws.Range("I66").Value = 542.5
$ws.Range("K66").Value = 4882.5
$ws.Range("M66").Value = -1138.5
$ws.Range("H92").Value = 203.2
$ws.Range("I92").Value = 209
$ws.Range("K92").Value = 627
$ws.Range("M92").Value = 621
$ws.Range("H113").Value = 1201.4286
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1201.4286
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3604.2858
$ws.Range("N113").Value = -7944.2858
$ws.Range("H121").Value = 1847.6428
$ws.Range("I121").Value = 1589.75
$ws.Range("J121").Value = 1890.625
$ws.Range("K121").Value = 4769.25
$ws.Range("L121").Value = 5671.875
$ws.Range("M121").Value = -3459.25
$ws.Range("N121").Value = -8291.875
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 500933.75
$ws.Range("I97").Value = 714579.3
$ws.Range("J97").Value = 2427.5
$ws.Range("K97").Value = 714579.3
$ws.Range("L97").Value = 2427.5
$ws.Range("M97").Value = -714083.3
$ws.Range("N97").Value = -3419.5
$ws.Range("H102").Value = 964.09375
$ws.Range("I102").Value = 824.2593000000001
$ws.Range("K102").Value = 824.2593000000001
$ws.Range("M102").Value = 797.7406999999999
$ws.Range("H126").Value = 3050.6365
$ws.Range("I126").Value = 2232
$ws.Range("J126").Value = 4033
$ws.Range("K126").Value = 6696
$ws.Range("L126").Value = 12099
$ws.Range("M126").Value = -4226
$ws.Range("N126").Value = -17039

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 47305
$ws.Range("I2").Value = 14000
$ws.Range("K2").Value = 14000
$ws.Range("M2").Value = -13888
$ws.Range("H22").Value = 1161.2667
$ws.Range("I22").Value = 944.5
$ws.Range("J22").Value = 1305.7778
$ws.Range("K22").Value = 944.5
$ws.Range("L22").Value = 1305.7778
$ws.Range("M22").Value = -649.5
$ws.Range("N22").Value = -1895.7778
$ws.Range("H27").Value = 1161.2667
$ws.Range("I27").Value = 944.5
$ws.Range("J27").Value = 1305.7778
$ws.Range("K27").Value = 944.5
$ws.Range("L27").Value = 1305.7778
$ws.Range("M27").Value = -837.5
$ws.Range("N27").Value = -1519.7778
$ws.Range("H93").Value = 2017.5555
$ws.Range("I93").Value = 2194.75
$ws.Range("J93").Value = 600
$ws.Range("K93").Value = 2194.75
$ws.Range("L93").Value = 600
$ws.Range("M93").Value = -946.75
$ws.Range("N93").Value = -3096

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 963.125
$ws.Range("I2").Value = 963.125
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 963.125
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -851.125
$ws.Range("H107").Value = 7986.7383
$ws.Range("I107").Value = 7925
$ws.Range("K107").Value = 23775
$ws.Range("M107").Value = -21855
$ws.Range("H132").Value = 2175565.5
$ws.Range("I132").Value = 1461.5555
$ws.Range("J132").Value = 3954377.8
$ws.Range("K132").Value = 4384.666499999999
$ws.Range("L132").Value = 11863133.4
$ws.Range("M132").Value = -1854.666499999999
$ws.Range("N132").Value = -11868193.4
$ws.Range("N2").ClearContents()
